$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "61.533.69"; ForceText = $False }
    @{ Cell = "E2"; Value = "  +0.52%  "; ForceText = $False }
    @{ Cell = "D3"; Value = "3.445.84"; ForceText = $False }
    @{ Cell = "E3"; Value = "  +1.30%  "; ForceText = $False }
    @{ Cell = "D4"; Value = "1.00"; ForceText = $True }
    @{ Cell = "E4"; Value = "  +0.00%  "; ForceText = $False }
    @{ Cell = "D5"; Value = "577.83"; ForceText = $True }
    @{ Cell = "E5"; Value = "  +0.93%  "; ForceText = $False }
    @{ Cell = "D6"; Value = "145.44"; ForceText = $True }
    @{ Cell = "E6"; Value = "  +4.68%  "; ForceText = $False }
    @{ Cell = "D7"; Value = "3.447.17"; ForceText = $False }
    @{ Cell = "E7"; Value = "  +1.38%  "; ForceText = $False }
    @{ Cell = "E8"; Value = "  +0.03%  "; ForceText = $False }
    @{ Cell = "E9"; Value = "  +2.45%  "; ForceText = $False }
    @{ Cell = "E10"; Value = "  +0.01%  "; ForceText = $False }
    @{ Cell = "E11"; Value = "  +3.63%  "; ForceText = $False }
    @{ Cell = "D12"; Value = "0.389"; ForceText = $True }
    @{ Cell = "E12"; Value = "  +2.44%  "; ForceText = $False }
    @{ Cell = "D13"; Value = "4.033.71"; ForceText = $False }
    @{ Cell = "E13"; Value = "  +1.27%  "; ForceText = $False }
    @{ Cell = "D14"; Value = "28.44"; ForceText = $True }
    @{ Cell = "E14"; Value = "  +6.33%  "; ForceText = $False }
    @{ Cell = "E15"; Value = "  -0.38%  "; ForceText = $False }
    @{ Cell = "E16"; Value = "  +1.05%  "; ForceText = $False }
    @{ Cell = "D17"; Value = "3.449.95"; ForceText = $False }
    @{ Cell = "E17"; Value = "  +1.34%  "; ForceText = $False }
    @{ Cell = "D18"; Value = "61.685.56"; ForceText = $False }
    @{ Cell = "E18"; Value = "  +0.70%  "; ForceText = $False }
    @{ Cell = "D19"; Value = "6.34"; ForceText = $True }
    @{ Cell = "E19"; Value = "  +6.77%  "; ForceText = $False }
    @{ Cell = "E20"; Value = "  +3.49%  "; ForceText = $False }
    @{ Cell = "D21"; Value = "9.44"; ForceText = $True }
    @{ Cell = "E21"; Value = "  +0.77%  "; ForceText = $False }
    @{ Cell = "D22"; Value = "403.17"; ForceText = $True }
    @{ Cell = "E22"; Value = "  +6.96%  "; ForceText = $False }
    @{ Cell = "D23"; Value = "0.569"; ForceText = $True }
    @{ Cell = "E23"; Value = "  +3.07%  "; ForceText = $False }
    @{ Cell = "D24"; Value = "74.45"; ForceText = $True }
    @{ Cell = "E24"; Value = "  +4.67%  "; ForceText = $False }
    @{ Cell = "E25"; Value = "  +0.09%  "; ForceText = $False }
    @{ Cell = "E26"; Value = "  +0.62%  "; ForceText = $False }
    @{ Cell = "E27"; Value = "  -0.17%  "; ForceText = $False }
    @{ Cell = "D28"; Value = "3.587.65"; ForceText = $False }
    @{ Cell = "E28"; Value = "  +1.62%  "; ForceText = $False }
    @{ Cell = "E29"; Value = "  +4.61%  "; ForceText = $False }
    @{ Cell = "D30"; Value = "7.63"; ForceText = $True }
    @{ Cell = "E30"; Value = "  +2.77%  "; ForceText = $False }
    @{ Cell = "D31"; Value = "1.00"; ForceText = $True }
    @{ Cell = "E31"; Value = "  +0.04%  "; ForceText = $False }
    @{ Cell = "D32"; Value = "8.27"; ForceText = $True }
    @{ Cell = "E32"; Value = "  +2.02%  "; ForceText = $False }
    @{ Cell = "E33"; Value = "  +1.96%  "; ForceText = $False }
    @{ Cell = "D34"; Value = "1.46"; ForceText = $True }
    @{ Cell = "E34"; Value = "  -9.97%  "; ForceText = $False }
    @{ Cell = "E35"; Value = "  -0.06%  "; ForceText = $False }
    @{ Cell = "D36"; Value = "23.90"; ForceText = $True }
    @{ Cell = "E36"; Value = "  +1.95%  "; ForceText = $False }
    @{ Cell = "D37"; Value = "7.05"; ForceText = $True }
    @{ Cell = "E37"; Value = "  +2.64%  "; ForceText = $False }
    @{ Cell = "D38"; Value = "3.474.08"; ForceText = $False }
    @{ Cell = "E38"; Value = "  +1.53%  "; ForceText = $False }
    @{ Cell = "B39"; Value = "ImmutableX"; ForceText = $False }
    @{ Cell = "C39"; Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; ForceText = $False }
    @{ Cell = "D39"; Value = "1.56"; ForceText = $True }
    @{ Cell = "E39"; Value = "  +0.17%  "; ForceText = $False }
    @{ Cell = "B40"; Value = "NEARProtocol"; ForceText = $False }
    @{ Cell = "C40"; Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"; ForceText = $False }
    @{ Cell = "D40"; Value = "5.14"; ForceText = $True }
    @{ Cell = "E40"; Value = "  +0.40%  "; ForceText = $False }
    @{ Cell = "D41"; Value = "167.22"; ForceText = $True }
    @{ Cell = "E41"; Value = "  +0.59%  "; ForceText = $False }
    @{ Cell = "D42"; Value = "0.0791"; ForceText = $True }
    @{ Cell = "E42"; Value = "  +2.77%  "; ForceText = $False }
    @{ Cell = "D43"; Value = "27.18"; ForceText = $True }
    @{ Cell = "E43"; Value = "  +4.26%  "; ForceText = $False }
    @{ Cell = "E44"; Value = "  +3.22%  "; ForceText = $False }
    @{ Cell = "E45"; Value = "  +3.14%  "; ForceText = $False }
    @{ Cell = "E46"; Value = "  -1.25%  "; ForceText = $False }
    @{ Cell = "E47"; Value = "  +0.03%  "; ForceText = $False }
    @{ Cell = "D48"; Value = "42.39"; ForceText = $True }
    @{ Cell = "E48"; Value = "  +1.02%  "; ForceText = $False }
    @{ Cell = "D49"; Value = "2.612.04"; ForceText = $False }
    @{ Cell = "E49"; Value = "  +3.45%  "; ForceText = $False }
    @{ Cell = "E50"; Value = "  -2.29%  "; ForceText = $False }
    @{ Cell = "E51"; Value = "  +2.51%  "; ForceText = $False }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $u.Value
}

